# Re-interpolated tree-diagram shape coordinates on slide 3 (see commit
# "corrected the interpolation of tree coords"). Slide indices below use
# 1-based PowerPoint COM numbering.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

$sh = $s.Shapes.Item("Connecteur droit 11")
$sh.Left = 57.618308
$sh.Top = 52.904371
$sh.Width = 589.129568
$sh.Height = 58.539174
$sh.VerticalFlip = 0

$sh = $s.Shapes.Item("Connecteur droit 13")
$sh.Left = 607.655631
$sh.Top = 111.443505
$sh.Width = 39.092245
$sh.Height = 364.933898

$sh = $s.Shapes.Item("Connecteur droit 15")
$sh.Left = 92.378780
$sh.Top = 127.848544
$sh.Width = 559.806654
$sh.Height = 5.567914

$sh = $s.Shapes.Item("Connecteur droit 17")
$sh.Left = 646.084371
$sh.Top = 83.745709
$sh.Width = 1.585631
$sh.Height = 390.570827

$sh = $s.Shapes.Item("ZoneTexte 24")
$sh.Left = 591.352717
$sh.Top = 470.611772
$sh.Width = 84.100827
$sh.Height = 23.628544

$sh = $s.Shapes.Item("ZoneTexte 25")
$sh.Left = 645.103505
$sh.Top = 390.349804
$sh.Width = 68.690512
$sh.Height = 23.628544

$sh = $s.Shapes.Item("ZoneTexte 26")
$sh.Left = 65.082481
$sh.Top = 98.331300
$sh.Width = 68.690512
$sh.Height = 23.628544

$sh = $s.Shapes.Item("Connecteur droit 6")
$sh.Left = 109.439961
$sh.Top = 416.286024
$sh.Width = 499.879961
$sh.Height = 58.030512

$sh = $s.Shapes.Item("Connecteur droit 9")
$sh.Left = 58.540591
$sh.Top = 52.904371
$sh.Width = 50.899410
$sh.Height = 363.426103

$sh = $s.Shapes.Item("Connecteur droit 19")
$sh.Left = 94.478859
$sh.Top = 474.316497
$sh.Width = 510.219725
$sh.Height = 0.000040
$sh.VerticalFlip = 0

$sh = $s.Shapes.Item("Connecteur droit 28")
$sh.Left = 607.655631
$sh.Top = 394.911930
$sh.Width = 0.792875
$sh.Height = 79.474607

$sh = $s.Shapes.Item("ZoneTexte 31")
$sh.Left = 544.535079
$sh.Top = 412.420040
$sh.Width = 68.690512
$sh.Height = 23.628544

$sh = $s.Shapes.Item("Connecteur droit avec flèche 34")
$sh.Left = 94.478859
$sh.Top = 301.768859
$sh.Width = 523.822009
$sh.Height = 75.542875

$sh = $s.Shapes.Item("ZoneTexte 37")
$sh.Left = 540.629489
$sh.Top = 367.759646
$sh.Width = 68.690512
$sh.Height = 23.628544

$sh = $s.Shapes.Item("Parallélogramme 1")
$sh.Left = 82.283505
$sh.Top = 24.291221
$sh.Width = 499.784843
$sh.Height = 480.699410

$sh = $s.Shapes.Item("Connecteur droit avec flèche 4")
$sh.Left = 94.478859
$sh.Top = 301.606024
$sh.Width = 499.629095
$sh.Height = 71.251930

$sh = $s.Shapes.Item("ZoneTexte 10")
$sh.Left = 527.728229
$sh.Top = 339.361694
$sh.Width = 68.690512
$sh.Height = 23.628544
